# Template Management & Send Message Tab Assertion update
# Update the sample "Ticket Id" values used in row 2 and row 3 of the
# Bulk Upload Template so they reflect new ticket ids.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell A2: 210321000128 -> 210321000131
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "210321000131"
$ws.Range("A2").Style = "Normal"

# Cell A3: 210321000129 -> 210321000132
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "210321000132"
$ws.Range("A3").Style = "Normal"
